$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.410.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.037.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.08"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.397"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0814"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.42"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +8.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.866"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.336.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.54"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.033.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.318.29"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.91"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0872"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.78"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.98%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.56"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.94"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.98%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0675"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.54%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.55"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.72%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.50"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0986"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +9.05%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0216"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.18%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.392.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +20.09%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.58"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.76%  "
